$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shapes are repositioned in EMUs. Shape.Left/.Top are expressed in points
# (1 pt = 12700 EMU) and stored internally as single-precision floats, so a
# naive EMU/12700 division can truncate one EMU short after the round-trip.
# Nudging by +0.5 EMU before the division keeps the saved value exact.
function Set-ShapePositionEmu {
    param($Shape, [double]$XEmu, [double]$YEmu)
    $Shape.Left = ($XEmu + 0.5) / 12700
    $Shape.Top  = ($YEmu + 0.5) / 12700
}

# --- Shape 1: "usuarios" -> "usuarioapp__c" (Salesforce object UsuarioApp__c) ---
$sh1 = $s.Shapes.Item(1)
Set-ShapePositionEmu $sh1 467544 392284
$tr1 = $sh1.TextFrame.TextRange
$tr1.Text = "u"
$tr1.InsertAfter("suarioapp__c") | Out-Null

# --- Shape 2: "plantas" -> "planta__c" (Salesforce object Planta__c) ---
$sh2 = $s.Shapes.Item(2)
Set-ShapePositionEmu $sh2 4932040 484620
$tr2 = $sh2.TextFrame.TextRange
$tr2.Text = "p"
$tr2.InsertAfter("lanta__c") | Out-Null

# --- Shape 3: "rutinas" -> "rutinas__c" (Salesforce object rutinas__c) ---
$sh3 = $s.Shapes.Item(3)
Set-ShapePositionEmu $sh3 4889188 2420888
$sh3.TextFrame.TextRange.Text = "rutinas__c"

# --- Shape 4: "tickets" - reposition only ---
$sh4 = $s.Shapes.Item(4)
Set-ShapePositionEmu $sh4 6012160 5301208

# --- Shape 5: "check-in-out-logs" - reposition only ---
$sh5 = $s.Shapes.Item(5)
Set-ShapePositionEmu $sh5 6012160 5949280
